$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "43.452.29"
Set-TextCell 2 5 "  +2.76%  "
Set-TextCell 3 4 "2.312.24"
Set-TextCell 3 5 "  +1.72%  "
Set-TextCell 4 5 "  +0.02%  "
Set-TextCell 5 4 "311.08"
Set-TextCell 5 5 "  +0.72%  "
Set-TextCell 6 4 "104.38"
Set-TextCell 6 5 "  +7.39%  "
Set-TextCell 8 5 "  +0.04%  "
Set-TextCell 9 5 "  +8.44%  "
Set-TextCell 10 4 "36.82"
Set-TextCell 10 5 "  +4.50%  "
Set-TextCell 11 4 "52.91"
Set-TextCell 11 5 "  +1.27%  "
Set-TextCell 12 5 "  +0.83%  "
Set-TextCell 13 5 "  -1.23%  "
Set-TextCell 14 4 "7.02"
Set-TextCell 14 5 "  +2.80%  "
Set-TextCell 15 4 "2.669.37"
Set-TextCell 15 5 "  +1.71%  "
Set-TextCell 16 4 "15.12"
Set-TextCell 16 5 "  +3.32%  "
Set-TextCell 17 4 "2.312.57"
Set-TextCell 17 5 "  +2.28%  "
Set-TextCell 18 4 "0.811"
Set-TextCell 18 5 "  +2.52%  "
Set-TextCell 19 4 "43.348.88"
Set-TextCell 19 5 "  +2.83%  "
Set-TextCell 20 4 "12.20"
Set-TextCell 20 5 "  -0.53%  "
Set-TextCell 21 5 "  +2.33%  "
Set-TextCell 22 4 "6.19"
Set-TextCell 22 5 "  +3.43%  "
Set-TextCell 23 4 "68.16"
Set-TextCell 23 5 "  +0.75%  "
Set-TextCell 24 4 "242.85"
Set-TextCell 24 5 "  +2.55%  "
Set-TextCell 25 5 "  +2.74%  "
Set-TextCell 26 4 "2.61"
Set-TextCell 26 5 "  +0.59%  "
Set-TextCell 27 5 "  +0.18%  "
Set-TextCell 28 2 "EthereumClassic"
Set-TextCell 28 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell 28 4 "24.91"
Set-TextCell 28 5 "  +5.60%  "
Set-TextCell 29 2 "Toncoin"
Set-TextCell 29 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell 29 4 "2.31"
Set-TextCell 29 5 "  +8.47%  "
Set-TextCell 30 2 "InjectiveProtocol"
Set-TextCell 30 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell 30 4 "37.12"
Set-TextCell 30 5 "  -0.43%  "
Set-TextCell 31 2 "Cosmos"
Set-TextCell 31 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 31 4 "9.64"
Set-TextCell 31 5 "  +0.86%  "
Set-TextCell 32 2 "Monero"
Set-TextCell 32 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 32 4 "167.75"
Set-TextCell 32 5 "  +2.34%  "
Set-TextCell 33 2 "Filecoin"
Set-TextCell 33 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 33 4 "5.29"
Set-TextCell 33 5 "  +0.69%  "
Set-TextCell 34 2 "FirstDigitalUSD"
Set-TextCell 34 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell 34 4 "1.00"
Set-TextCell 34 5 "  -0.06%  "
Set-TextCell 35 4 "18.40"
Set-TextCell 35 5 "  +4.07%  "
Set-TextCell 36 2 "WEMIXToken"
Set-TextCell 36 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell 36 4 "2.54"
Set-TextCell 36 5 "  +6.78%  "
Set-TextCell 37 2 "Hedera"
Set-TextCell 37 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 37 4 "0.0744"
Set-TextCell 37 5 "  +1.22%  "
Set-TextCell 38 2 "LidoDAOToken"
Set-TextCell 38 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell 38 4 "3.06"
Set-TextCell 38 5 "  -1.10%  "
Set-TextCell 39 2 "ARBITRUM"
Set-TextCell 39 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 39 4 "1.88"
Set-TextCell 39 5 "  +3.28%  "
Set-TextCell 40 2 "Kaspa"
Set-TextCell 40 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell 40 4 "0.106"
Set-TextCell 40 5 "  +2.29%  "
Set-TextCell 41 2 "RenderToken"
Set-TextCell 41 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 41 4 "4.46"
Set-TextCell 41 5 "  +6.49%  "
Set-TextCell 42 2 "Stellar"
Set-TextCell 42 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 42 4 "0.116"
Set-TextCell 42 5 "  +0.73%  "
Set-TextCell 43 2 "ApeXProtocol"
Set-TextCell 43 3 "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell 43 4 "2.72"
Set-TextCell 43 5 "  +19.70%  "
Set-TextCell 44 2 "VeChain"
Set-TextCell 44 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 44 4 "0.0293"
Set-TextCell 44 5 "  +3.65%  "
Set-TextCell 45 2 "Maker"
Set-TextCell 45 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell 45 4 "1.990.31"
Set-TextCell 45 5 "  +2.13%  "
Set-TextCell 46 2 "EnergySwap"
Set-TextCell 46 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell 46 4 "18.99"
Set-TextCell 46 5 "  +0.95%  "
Set-TextCell 47 2 "NEARProtocol"
Set-TextCell 47 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 47 4 "3.08"
Set-TextCell 47 5 "  +4.00%  "
Set-TextCell 48 2 "FraxShare"
Set-TextCell 48 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 48 4 "10.01"
Set-TextCell 48 5 "  +2.25%  "
Set-TextCell 49 2 "MultiversX"
Set-TextCell 49 3 "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell 49 4 "55.89"
Set-TextCell 49 5 "  +3.87%  "
Set-TextCell 50 2 "HuobiToken"
Set-TextCell 50 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 50 4 "2.95"
Set-TextCell 50 5 "  +2.13%  "
Set-TextCell 51 2 "Stacks"
Set-TextCell 51 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell 51 4 "1.60"
Set-TextCell 51 5 "  +8.93%  "
